$wb = $excel.ActiveWorkbook

# Sheet "展览": update "想去人数" (F column) figures
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2881
$ws1.Range("F4").Value = 132
$ws1.Range("F5").Value = 49

# Sheet "全部类型": update the same rows (duplicated data) in the combined sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 2881
$ws4.Range("F8").Value = 132
$ws4.Range("F10").Value = 49
